# Applies the Jan 2024 update to fl2_meta.xlsx:
#  - renames the worksheet from "Sheet1" to "fl2_meta"
#  - bumps the "Treatment" (column B) values for rows 26/34/42/50/58
#  - corrects A78 from the numeric 70 to the text label "70B"
#  - updates the saved selection/scroll position to H76 (no more scrolled view)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "fl2_meta"

# Treatment column corrections
$ws.Range("B26").Value = 2
$ws.Range("B34").Value = 3
$ws.Range("B42").Value = 4
$ws.Range("B50").Value = 5
$ws.Range("B58").Value = 6

# Syringe label correction: 70 -> 70B
$ws.Range("A78").Value = "70B"

# Update the active cell / selection (also resets the scrolled top-left cell)
$ws.Range("H76").Select() | Out-Null
